$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 1).Value = 111576450
$ws.Cells.Item(13, 2).Value = 96348
$ws.Cells.Item(13, 4).Value = "VU"
$ws.Cells.Item(13, 5).Value = 220787
$ws.Cells.Item(13, 6).Value = "Knärot"
$ws.Cells.Item(13, 7).Value = "Goodyera repens"
$ws.Cells.Item(13, 8).Value = "(L.) R. Br."
$ws.Cells.Item(13, 17).Value = 562979.5212303887
$ws.Cells.Item(13, 18).Value = 6954739.97881452
$ws.Cells.Item(13, 26).Value = "17:10"
$ws.Cells.Item(13, 28).Value = "17:10"
$ws.Cells.Item(13, 29).Value = "Rikligt"
$ws.Cells.Item(14, 1).Value = 111573569
$ws.Cells.Item(14, 2).Value = 96348
$ws.Cells.Item(14, 4).Value = "VU"
$ws.Cells.Item(14, 5).Value = 220787
$ws.Cells.Item(14, 6).Value = "Knärot"
$ws.Cells.Item(14, 7).Value = "Goodyera repens"
$ws.Cells.Item(14, 8).Value = "(L.) R. Br."
$ws.Cells.Item(14, 17).Value = 562701.9737813871
$ws.Cells.Item(14, 18).Value = 6954788.374143652
$ws.Cells.Item(14, 26).Value = "00:00"
$ws.Cells.Item(14, 28).Value = "00:00"
$ws.Cells.Item(15, 1).Value = 111574128
$ws.Cells.Item(15, 2).Value = 96348
$ws.Cells.Item(15, 4).Value = "VU"
$ws.Cells.Item(15, 5).Value = 220787
$ws.Cells.Item(15, 6).Value = "Knärot"
$ws.Cells.Item(15, 7).Value = "Goodyera repens"
$ws.Cells.Item(15, 8).Value = "(L.) R. Br."
$ws.Cells.Item(15, 9).ClearContents() | Out-Null
$ws.Cells.Item(15, 17).Value = 562555.4143375416
$ws.Cells.Item(15, 18).Value = 6954835.60431945
$ws.Cells.Item(15, 29).ClearContents() | Out-Null
$ws.Cells.Item(16, 1).Value = 111578127
$ws.Cells.Item(16, 2).Value = 56543
$ws.Cells.Item(16, 4).Value = "NT"
$ws.Cells.Item(16, 5).Value = 103021
$ws.Cells.Item(16, 6).Value = "Talltita"
$ws.Cells.Item(16, 7).Value = "Poecile montanus"
$ws.Cells.Item(16, 8).Value = "(Conrad von Baldenstein, 1827)"
$ws.Cells.Item(16, 13).Value = "lockläte, övriga läten"
$ws.Cells.Item(16, 17).Value = 562937.8205991766
$ws.Cells.Item(16, 18).Value = 6954541.406048392
$ws.Cells.Item(16, 26).Value = "18:30"
$ws.Cells.Item(16, 28).Value = "18:30"
$ws.Cells.Item(17, 1).Value = 111575796
$ws.Cells.Item(17, 17).Value = 562855.7640570825
$ws.Cells.Item(17, 18).Value = 6954651.349091855
$ws.Cells.Item(17, 26).Value = "16:39"
$ws.Cells.Item(17, 28).Value = "16:39"
$ws.Cells.Item(18, 1).Value = 111573866
$ws.Cells.Item(18, 2).Value = 96348
$ws.Cells.Item(18, 4).Value = "VU"
$ws.Cells.Item(18, 5).Value = 220787
$ws.Cells.Item(18, 6).Value = "Knärot"
$ws.Cells.Item(18, 7).Value = "Goodyera repens"
$ws.Cells.Item(18, 8).Value = "(L.) R. Br."
$ws.Cells.Item(18, 17).Value = 562601.7570288588
$ws.Cells.Item(18, 18).Value = 6954814.918206804
$ws.Cells.Item(18, 26).Value = "15:17"
$ws.Cells.Item(18, 28).Value = "15:17"
$ws.Cells.Item(19, 1).Value = 111578197
$ws.Cells.Item(19, 2).Value = 96348
$ws.Cells.Item(19, 4).Value = "VU"
$ws.Cells.Item(19, 5).Value = 220787
$ws.Cells.Item(19, 6).Value = "Knärot"
$ws.Cells.Item(19, 7).Value = "Goodyera repens"
$ws.Cells.Item(19, 8).Value = "(L.) R. Br."
$ws.Cells.Item(19, 17).Value = 563026.0554397166
$ws.Cells.Item(19, 18).Value = 6954541.256262898
$ws.Cells.Item(19, 26).Value = "00:00"
$ws.Cells.Item(19, 28).Value = "00:00"
$ws.Cells.Item(20, 1).Value = 111575785
$ws.Cells.Item(20, 2).Value = 89845
$ws.Cells.Item(20, 5).Value = 1209
$ws.Cells.Item(20, 6).Value = "Rynkskinn"
$ws.Cells.Item(20, 7).Value = "Phlebia centrifuga"
$ws.Cells.Item(20, 8).Value = "P.Karst."
$ws.Cells.Item(20, 17).Value = 562859.2727272335
$ws.Cells.Item(20, 18).Value = 6954660.134623887
$ws.Cells.Item(20, 26).Value = "16:39"
$ws.Cells.Item(20, 28).Value = "16:39"
$ws.Cells.Item(21, 1).Value = 111573803
$ws.Cells.Item(21, 17).Value = 562591.0245237258
$ws.Cells.Item(21, 18).Value = 6954847.751526525
$ws.Cells.Item(21, 26).Value = "15:14"
$ws.Cells.Item(21, 28).Value = "15:14"
$ws.Cells.Item(22, 1).Value = 111574509
$ws.Cells.Item(22, 2).Value = 96348
$ws.Cells.Item(22, 4).Value = "VU"
$ws.Cells.Item(22, 5).Value = 220787
$ws.Cells.Item(22, 6).Value = "Knärot"
$ws.Cells.Item(22, 7).Value = "Goodyera repens"
$ws.Cells.Item(22, 8).Value = "(L.) R. Br."
$ws.Cells.Item(22, 17).Value = 562529.1073683554
$ws.Cells.Item(22, 18).Value = 6954769.030357216
$ws.Cells.Item(22, 26).Value = "15:45"
$ws.Cells.Item(22, 28).Value = "15:45"
$ws.Cells.Item(23, 1).Value = 111575868
$ws.Cells.Item(23, 17).Value = 562854.9195222461
$ws.Cells.Item(23, 18).Value = 6954623.341454657
$ws.Cells.Item(23, 26).Value = "16:43"
$ws.Cells.Item(23, 28).Value = "16:43"
$ws.Cells.Item(24, 1).Value = 111574429
$ws.Cells.Item(24, 2).Value = 89405
$ws.Cells.Item(24, 4).Value = "NT"
$ws.Cells.Item(24, 5).Value = 1202
$ws.Cells.Item(24, 6).Value = "Ullticka"
$ws.Cells.Item(24, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(24, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(24, 17).Value = 562547.0565141424
$ws.Cells.Item(24, 18).Value = 6954767.535469687
$ws.Cells.Item(24, 26).Value = "15:42"
$ws.Cells.Item(24, 28).Value = "15:42"
$ws.Cells.Item(25, 1).Value = 111574334
$ws.Cells.Item(25, 2).Value = 89405
$ws.Cells.Item(25, 4).Value = "NT"
$ws.Cells.Item(25, 5).Value = 1202
$ws.Cells.Item(25, 6).Value = "Ullticka"
$ws.Cells.Item(25, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(25, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(25, 17).Value = 562557.3535548041
$ws.Cells.Item(25, 18).Value = 6954757.635990249
$ws.Cells.Item(25, 26).Value = "15:26"
$ws.Cells.Item(25, 28).Value = "15:26"
$ws.Cells.Item(26, 1).Value = 111574403
$ws.Cells.Item(26, 2).Value = 89686
$ws.Cells.Item(26, 4).Value = "NT"
$ws.Cells.Item(26, 5).Value = 658
$ws.Cells.Item(26, 6).Value = "Rosenticka"
$ws.Cells.Item(26, 7).Value = "Rhodofomes roseus"
$ws.Cells.Item(26, 8).Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Cells.Item(26, 17).Value = 562547.0565141424
$ws.Cells.Item(26, 18).Value = 6954767.535469687
$ws.Cells.Item(26, 26).Value = "15:26"
$ws.Cells.Item(26, 28).Value = "15:26"
$ws.Cells.Item(27, 1).Value = 111573948
$ws.Cells.Item(27, 17).Value = 562576.2301468613
$ws.Cells.Item(27, 18).Value = 6954852.517936011
$ws.Cells.Item(27, 26).Value = "15:20"
$ws.Cells.Item(27, 28).Value = "15:20"
$ws.Cells.Item(28, 1).Value = 111576771
$ws.Cells.Item(28, 17).Value = 562807.4867926922
$ws.Cells.Item(28, 18).Value = 6954821.585021482
$ws.Cells.Item(28, 26).Value = "17:24"
$ws.Cells.Item(28, 28).Value = "17:24"
$ws.Cells.Item(29, 1).Value = 111576401
$ws.Cells.Item(29, 2).Value = 89369
$ws.Cells.Item(29, 4).Value = "LC"
$ws.Cells.Item(29, 5).Value = 5447
$ws.Cells.Item(29, 6).Value = "Vedticka"
$ws.Cells.Item(29, 7).Value = "Fuscoporia viticola"
$ws.Cells.Item(29, 8).Value = "(Schwein.) Murrill"
$ws.Cells.Item(29, 17).Value = 562964.914807545
$ws.Cells.Item(29, 18).Value = 6954710.791209211
$ws.Cells.Item(29, 26).Value = "16:51"
$ws.Cells.Item(29, 28).Value = "16:51"
$ws.Cells.Item(30, 1).Value = 111576037
$ws.Cells.Item(30, 2).Value = 89686
$ws.Cells.Item(30, 4).Value = "NT"
$ws.Cells.Item(30, 5).Value = 658
$ws.Cells.Item(30, 6).Value = "Rosenticka"
$ws.Cells.Item(30, 7).Value = "Rhodofomes roseus"
$ws.Cells.Item(30, 8).Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Cells.Item(30, 17).Value = 562852.9463231879
$ws.Cells.Item(30, 18).Value = 6954606.325244571
$ws.Cells.Item(30, 26).Value = "16:51"
$ws.Cells.Item(30, 28).Value = "16:51"
$ws.Cells.Item(31, 1).Value = 111574338
$ws.Cells.Item(31, 2).Value = 89686
$ws.Cells.Item(31, 5).Value = 658
$ws.Cells.Item(31, 6).Value = "Rosenticka"
$ws.Cells.Item(31, 7).Value = "Rhodofomes roseus"
$ws.Cells.Item(31, 8).Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Cells.Item(31, 13).ClearContents() | Out-Null
$ws.Cells.Item(31, 17).Value = 562557.3535548041
$ws.Cells.Item(31, 18).Value = 6954757.635990249
$ws.Cells.Item(31, 26).Value = "15:26"
$ws.Cells.Item(31, 28).Value = "15:26"
$ws.Cells.Item(32, 1).Value = 111574240
$ws.Cells.Item(32, 2).Value = 56543
$ws.Cells.Item(32, 4).Value = "NT"
$ws.Cells.Item(32, 5).Value = 103021
$ws.Cells.Item(32, 6).Value = "Talltita"
$ws.Cells.Item(32, 7).Value = "Poecile montanus"
$ws.Cells.Item(32, 8).Value = "(Conrad von Baldenstein, 1827)"
$ws.Cells.Item(32, 9).Value = "5"
$ws.Cells.Item(32, 17).Value = 562533.1227179464
$ws.Cells.Item(32, 18).Value = 6954848.029061474
$ws.Cells.Item(32, 26).Value = "15:26"
$ws.Cells.Item(32, 28).Value = "15:26"
$ws.Cells.Item(32, 29).Value = "Familj med 5 talltitor. Permanent revir"
$ws.Cells.Item(33, 1).Value = 111574689
$ws.Cells.Item(33, 17).Value = 562517.0252856832
$ws.Cells.Item(33, 18).Value = 6954776.14289257
$ws.Cells.Item(33, 26).Value = "15:47"
$ws.Cells.Item(33, 28).Value = "15:47"
